$wb = $excel.ActiveWorkbook

# --- Add the new worksheet "User - Pengaturan Email" after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "User - Pengaturan Email"

# --- Header row ---
$ws.Range("A1").Value = "email"
$ws.Range("B1").Value = "kondisi"
$ws.Range("C1").Value = "keterangan"

# --- Row 2: empty email ---
$ws.Range("B2").Value = "fail"
$ws.Range("C2").Value = "emptyEmail"

# --- Row 3: invalid email ---
$ws.Range("A3").Value = "yoke"
$ws.Range("B3").Value = "fail"
$ws.Range("C3").Value = "invalidEmail"

# --- Row 4: cancel ---
$ws.Range("A4").Value = "dennaleksanti@gmail.com"
$ws.Range("B4").Value = "fail"
$ws.Range("C4").Value = "cancel"

# --- Row 5: pass ---
$ws.Range("A5").Value = "yokebethdenna@gmail.com"
$ws.Range("B5").Value = "pass"

# --- Hyperlinks for the two email addresses ---
$null = $ws.Hyperlinks.Add($ws.Range("A4"), "mailto:dennaleksanti@gmail.com")
$null = $ws.Hyperlinks.Add($ws.Range("A5"), "mailto:yokebethdenna@gmail.com")

# --- Column width for column A ---
$ws.Columns.Item(1).ColumnWidth = 27.109375

# --- Final selection on this sheet (also makes it the active/visible sheet) ---
$null = $ws.Range("C5").Select()

Write-Host "Added sheet 'User - Pengaturan Email'"
